$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" row (row 2, column B) holds the Cypher query used to build the
# Cases report. This edit drops the trailing "Cohort" column (and the OPTIONAL
# MATCH's co binding usage in the RETURN) from that query, per the commit:
# "Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study".

$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
            "MATCH (c)<--(diag:diagnosis)`n" +
            "MATCH (samp:sample)-->(c) `n" +
            " WHERE samp.specific_sample_pathology IN [`"Astrocytoma`"]  `n" +
            "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
            "  WITH DISTINCT c, s, demo, diag, co`n" +
            "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
            "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
            "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
            "        coalesce(demo.breed, '') AS Breed ,`n" +
            "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
            "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
            "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
            "        coalesce(demo.sex, '') AS Sex ,`n" +
            "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
            "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
            "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newQuery

$ws.Range("B2").Select()
